$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MaddisonJourney")

$ws.Range("B4").Value  = "cum with me right now Daddy... this is just for you 💕"
$ws.Range("B6").Value  = "I'm about to cum so hard Daddy... watch me, please don't look away"
$ws.Range("B7").Value  = "don't leave... I'm so close and I need to feel you right here when I finish 💕"
$ws.Range("B8").Value  = "oh fuck oh fuck 🥺"
$ws.Range("B9").Value  = "I've never gone this far with anyone Daddy... watch what you made me do 💕"
$ws.Range("B11").Value = "I'm almost there and I need you to see what you're doing to me right now"
$ws.Range("B12").Value = "I'm going faster and faster and I can feel myself getting closer... my whole body is trembling 💕"
$ws.Range("B13").Value = "my fingers are inside my pussy and I'm moaning so loud right now... I hope nobody can hear me"
$ws.Range("B14").Value = "fuckk 🥺"
$ws.Range("B15").Value = "see what you did to me? I can't stop 💕"
$ws.Range("B17").Value = "tell me exactly what you're thinking right now Daddy... I want to hear everything while I touch myself"
$ws.Range("B18").Value = "every part of me is on fire right now and it keeps getting more intense because of you 💕"
$ws.Range("B19").Value = "I can't stop now even if I wanted to... my hand is already between my thighs and I'm soaked"
$ws.Range("B20").Value = "omg... I can't believe that just happened 🥺"
$ws.Range("B21").Value = "this is what you're making me do... I can't believe I'm showing you this 💕"
$ws.Range("B23").Value = "I'm touching myself right now and I want you to know it's because of you Daddy"
$ws.Range("B24").Value = "my breathing is getting heavier and I keep arching my back... my body wants something and I think it's you"
$ws.Range("B25").Value = "wait you actually liked that? something just shifted inside me and I can feel it everywhere 🥺"
